$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Existing rows 2-51: column I (trade_date) switches number format
#    from "YYYY-MM-DD" (style s=3) to "YYYY-MM-DD HH:MM:SS" (style s=2),
#    matching the format already used by column B (datetime).
$ws.Range("I2:I51").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2) Append four new data rows (52-55) with the same column layout as the
#    rest of the sheet:
#    A close, B datetime, C exchange_code, D high, E low, F open,
#    G stock_code, H volume, I trade_date, J isin, K company, L symbol, M source
$newRows = @(
  @(443.9,  46049, "NSE", 454.3,  438.8,   451.1, "TATCOV", 12022746, 46049, "INE1TAE01010", "Tata Motors Ltd", "TATCOV", "BREEZE"),
  @(467.95, 46050, "NSE", 470.25, 448.3,   449.3, "TATCOV", 12586112, 46050, "INE1TAE01010", "Tata Motors Ltd", "TATCOV", "BREEZE"),
  @(470.2,  46051, "NSE", 475.4,  459.75,  473.5, "TATCOV", 11600368, 46051, "INE1TAE01010", "Tata Motors Ltd", "TATCOV", "BREEZE"),
  @(458.5,  46052, "NSE", 485.3,  451,     457,   "TATCOV", 15731508, 46052, "INE1TAE01010", "Tata Motors Ltd", "TATCOV", "BREEZE")
)

$startRow = 52
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $ws.Cells.Item($r, 8).Value2 = $row[7]
    $ws.Cells.Item($r, 9).Value2 = $row[8]
    $ws.Cells.Item($r, 10).Value2 = $row[9]
    $ws.Cells.Item($r, 11).Value2 = $row[10]
    $ws.Cells.Item($r, 12).Value2 = $row[11]
    $ws.Cells.Item($r, 13).Value2 = $row[12]

    # New rows keep the original per-column number formats:
    # B (datetime) = "YYYY-MM-DD HH:MM:SS" (style s=2)
    # I (trade_date) = "YYYY-MM-DD" (style s=3)
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 9).NumberFormat = "YYYY-MM-DD"
}
